$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# A leading apostrophe forces Excel to store a numeric-looking Price as text,
# matching the original inline-string cell type (no "#" -> number coercion).

$ws.Range("D2").Value = "42.551.34"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "2.298.57"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'322.86"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").Value = "'104.69"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'40.21"
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("D11").Value = "'0.0906"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Value = "'0.975"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "'15.32"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "2.646.94"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "2.287.91"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "42.662.79"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'13.45"
$ws.Range("E21").Value = "  +35.93%  "
$ws.Range("D22").Value = "'73.63"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'3.61"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'270.83"
$ws.Range("E24").Value = "  -4.02%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'22.62"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'38.32"
$ws.Range("E30").Value = "  +11.83%  "
$ws.Range("D31").Value = "'165.43"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  +6.29%  "
$ws.Range("D33").Value = "'0.0886"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -12.65%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "'4.62"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "'0.0356"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "'1.55"
$ws.Range("E41").Value = "  +6.71%  "
$ws.Range("D42").Value = "'98.58"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "'70.30"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  +4.88%  "
$ws.Range("D47").Value = "'82.45"
$ws.Range("E47").Value = "  +8.82%  "
$ws.Range("D48").Value = "'113.38"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'5.28"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").Value = "1.594.71"
$ws.Range("E51").Value = "  +4.12%  "
